$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'41.875.59"
$ws.Range('E2').Value = '  +4.81%  '
$ws.Range('D3').Value = "'2.267.77"
$ws.Range('E3').Value = '  +2.20%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'302.28"
$ws.Range('E5').Value = '  +3.44%  '
$ws.Range('D6').Value = "'92.31"
$ws.Range('E6').Value = '  +6.12%  '
$ws.Range('D7').Value = "'0.532"
$ws.Range('E7').Value = '  +3.49%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('E10').Value = '  +8.27%  '
$ws.Range('D11').Value = "'32.27"
$ws.Range('E11').Value = '  +5.84%  '
$ws.Range('E12').Value = '  +2.35%  '
$ws.Range('E13').Value = '  +2.05%  '
$ws.Range('D14').Value = "'6.67"
$ws.Range('E14').Value = '  +3.62%  '
$ws.Range('D15').Value = "'2.616.13"
$ws.Range('E15').Value = '  +2.06%  '
$ws.Range('D16').Value = "'14.18"
$ws.Range('E16').Value = '  +3.05%  '
$ws.Range('D17').Value = "'2.275.63"
$ws.Range('E17').Value = '  +4.69%  '
$ws.Range('D18').Value = "'0.759"
$ws.Range('E18').Value = '  +3.56%  '
$ws.Range('D19').Value = "'41.784.18"
$ws.Range('E19').Value = '  +4.80%  '
$ws.Range('D20').Value = "'12.09"
$ws.Range('E20').Value = '  +8.26%  '
$ws.Range('E21').Value = '  +1.94%  '
$ws.Range('D22').Value = "'5.96"
$ws.Range('E22').Value = '  +3.55%  '
$ws.Range('D23').Value = "'66.99"
$ws.Range('E23').Value = '  +2.24%  '
$ws.Range('D24').Value = "'241.97"
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').Value = "'2.55"
$ws.Range('E25').Value = '  +3.64%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  +4.19%  '
$ws.Range('D28').Value = "'23.90"
$ws.Range('E28').Value = '  +2.87%  '
$ws.Range('E29').Value = '  +4.38%  '
$ws.Range('E30').Value = '  -12.04%  '
$ws.Range('D31').Value = "'159.31"
$ws.Range('E31').Value = '  +1.09%  '
$ws.Range('E32').Value = '  +6.19%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').Value = "'5.18"
$ws.Range('E34').Value = '  +4.16%  '
$ws.Range('E35').Value = '  +4.40%  '
$ws.Range('E36').Value = '  +3.19%  '
$ws.Range('D37').Value = "'2.38"
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('E38').Value = '  +5.93%  '
$ws.Range('E39').Value = '  +3.62%  '
$ws.Range('D40').Value = "'16.54"
$ws.Range('E40').Value = '  +9.03%  '
$ws.Range('D41').Value = "'1.82"
$ws.Range('E41').Value = '  +4.89%  '
$ws.Range('D42').Value = "'3.92"
$ws.Range('E42').Value = '  +5.42%  '
$ws.Range('D43').Value = "'2.074.74"
$ws.Range('E43').Value = '  -0.73%  '
$ws.Range('D44').Value = "'19.64"
$ws.Range('E44').Value = '  +9.82%  '
$ws.Range('E45').Value = '  +3.45%  '
$ws.Range('E46').Value = '  +3.81%  '
$ws.Range('D47').Value = "'2.93"
$ws.Range('E47').Value = '  +8.52%  '
$ws.Range('D48').Value = "'2.04"
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('E49').Value = '  +3.58%  '
$ws.Range('E50').Value = '  +3.23%  '
$ws.Range('D51').Value = "'51.83"
$ws.Range('E51').Value = '  +5.71%  '
